$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''67.634.41'
$ws.Range("E2").Value = '''  -7.10%  '

# Row 3
$ws.Range("D3").Value = '''3.675.21'
$ws.Range("E3").Value = '''  -6.99%  '

# Row 4
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '''  +0.43%  '

# Row 5
$ws.Range("D5").Value = '''578.90'
$ws.Range("E5").Value = '''  -5.00%  '

# Row 6
$ws.Range("D6").Value = '''174.40'
$ws.Range("E6").Value = '''  +4.68%  '

# Row 7
$ws.Range("D7").Value = '''3.668.67'
$ws.Range("E7").Value = '''  -6.99%  '

# Row 8
$ws.Range("D8").Value = '''0.629'
$ws.Range("E8").Value = '''  -6.82%  '

# Row 9
$ws.Range("E9").Value = '''  +0.28%  '

# Row 10
$ws.Range("D10").Value = '''0.711'
$ws.Range("E10").Value = '''  -4.78%  '

# Row 11
$ws.Range("D11").Value = '''0.164'
$ws.Range("E11").Value = '''  -7.67%  '

# Row 12
$ws.Range("D12").Value = '''52.25'
$ws.Range("E12").Value = '''  -7.05%  '

# Row 13
$ws.Range("D13").Value = '''0.0000298'
$ws.Range("E13").Value = '''  -10.12%  '

# Row 14
$ws.Range("D14").Value = '''10.61'
$ws.Range("E14").Value = '''  -3.71%  '

# Row 15
$ws.Range("D15").Value = '''4.286.65'
$ws.Range("E15").Value = '''  -6.51%  '

# Row 16
$ws.Range("D16").Value = '''3.679.46'

# Row 17
$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").Value = '''0.127'
$ws.Range("E17").Value = '''  -2.91%  '

# Row 18
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").Value = '''19.33'
$ws.Range("E18").Value = '''  -5.10%  '

# Row 19
$ws.Range("D19").Value = '''12.96'
$ws.Range("E19").Value = '''  -7.05%  '

# Row 20
$ws.Range("D20").Value = '''1.13'
$ws.Range("E20").Value = '''  -8.36%  '

# Row 21
$ws.Range("D21").Value = '''67.617.11'
$ws.Range("E21").Value = '''  -6.95%  '

# Row 22
$ws.Range("D22").Value = '''407.25'
$ws.Range("E22").Value = '''  -6.37%  '

# Row 23
$ws.Range("D23").Value = '''4.53'
$ws.Range("E23").Value = '''  -6.45%  '

# Row 24
$ws.Range("D24").Value = '''88.01'
$ws.Range("E24").Value = '''  -7.41%  '

# Row 25
$ws.Range("D25").Value = '''3.06'
$ws.Range("E25").Value = '''  -8.65%  '

# Row 26
$ws.Range("D26").Value = '''12.76'
$ws.Range("E26").Value = '''  -9.58%  '

# Row 27
$ws.Range("D27").Value = '''10.71'
$ws.Range("E27").Value = '''  -3.61%  '

# Row 28
$ws.Range("E28").Value = '''  -6.51%  '

# Row 29
$ws.Range("D29").Value = '''5.96'
$ws.Range("E29").Value = '''  +0.22%  '

# Row 30
$ws.Range("D30").Value = '''9.50'
$ws.Range("E30").Value = '''  -8.48%  '

# Row 31
$ws.Range("D31").Value = '''8.13'
$ws.Range("E31").Value = '''  +2.11%  '

# Row 32
$ws.Range("D32").Value = '''32.80'
$ws.Range("E32").Value = '''  -8.29%  '

# Row 33
$ws.Range("D33").Value = '''12.67'
$ws.Range("E33").Value = '''  -6.74%  '

# Row 34
$ws.Range("B34").Value = 'InjectiveProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D34").Value = '''44.47'
$ws.Range("E34").Value = '''  -6.85%  '

# Row 35
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '''0.118'
$ws.Range("E35").Value = '''  -8.91%  '

# Row 36
$ws.Range("B36").Value = 'OKB'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D36").Value = '''65.73'
$ws.Range("E36").Value = '''  -5.95%  '

# Row 37
$ws.Range("D37").Value = '''0.0₃0928'
$ws.Range("E37").Value = '''  -7.62%  '

# Row 38
$ws.Range("D38").Value = '''587.90'
$ws.Range("E38").Value = '''  -7.73%  '

# Row 39
$ws.Range("D39").Value = '''0.402'
$ws.Range("E39").Value = '''  -6.52%  '

# Row 40
$ws.Range("D40").Value = '''1.00'
$ws.Range("E40").Value = '''  +0.00%  '

# Row 41
$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D41").Value = '''1.00'
$ws.Range("E41").Value = '''  +0.24%  '

# Row 42
$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").Value = '''3.23'
$ws.Range("E42").Value = '''  +13.13%  '

# Row 43
$ws.Range("D43").Value = '''0.136'
$ws.Range("E43").Value = '''  -7.07%  '

# Row 44
$ws.Range("D44").Value = '''3.05'
$ws.Range("E44").Value = '''  -11.85%  '

# Row 45
$ws.Range("D45").Value = '''0.0440'
$ws.Range("E45").Value = '''  -8.44%  '

# Row 46
$ws.Range("B46").Value = 'Fetch.AI'
$ws.Range("C46").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D46").Value = '''2.58'
$ws.Range("E46").Value = '''  -0.49%  '

# Row 47
$ws.Range("B47").Value = 'THORChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D47").Value = '''9.47'
$ws.Range("E47").Value = '''  -11.74%  '

# Row 48
$ws.Range("D48").Value = '''0.134'
$ws.Range("E48").Value = '''  -8.98%  '

# Row 49
$ws.Range("D49").Value = '''2.731.16'
$ws.Range("E49").Value = '''  -4.42%  '

# Row 50
$ws.Range("B50").Value = 'ApeXProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D50").Value = '''3.14'
$ws.Range("E50").Value = '''  -7.89%  '

# Row 51
$ws.Range("B51").Value = 'WEMIXToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D51").Value = '''2.66'
$ws.Range("E51").Value = '''  -16.07%  '
